$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note cells next to the "syringe Actuation" / PressureCalibration test row (row 32)
# and a new row 33 with an additional testing note.
$ws.Range("E32").Value = "(create another script to sycn RTC?)"
$ws.Range("E33").Value = "(false triggering?)"

# Update the active selection to match the new edit location
$ws.Range("E33").Select()
